$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows just above row 149 (existing rows 149..244 shift down to 151..246)
$ws.Rows.Item(149).Insert()
$ws.Rows.Item(149).Insert()

# --- Row 149: new weekly record (Paine, 1a nueva(o)) ---
$ws.Range("A149").Value = 4
$ws.Range("B149").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C149").Value = "Los Lagos"
$ws.Range("D149").Value = 44582
$ws.Range("E149").Value = 10
$ws.Range("F149").Value = 100112045
$ws.Range("G149").Value = "Zapallo"
$ws.Range("H149").Value = "Paine"
$ws.Range("I149").Value = "1a nueva(o)"
$ws.Range("J149").Value = 600
$ws.Range("K149").Value = 500
$ws.Range("L149").Value = 500
$ws.Range("M149").Value = 500
$ws.Range("N149").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O149").Value = "Región de O'Higgins"
$ws.Range("P149").Value = 500
$ws.Range("Q149").Value = 1
$ws.Range("R149").Value = "Hortaliza"

# --- Row 150: new weekly record (Paine, 2a nueva(o)) ---
$ws.Range("A150").Value = 4
$ws.Range("B150").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C150").Value = "Los Lagos"
$ws.Range("D150").Value = 44582
$ws.Range("E150").Value = 10
$ws.Range("F150").Value = 100112045
$ws.Range("G150").Value = "Zapallo"
$ws.Range("H150").Value = "Paine"
$ws.Range("I150").Value = "2a nueva(o)"
$ws.Range("J150").Value = 600
$ws.Range("K150").Value = 350
$ws.Range("L150").Value = 350
$ws.Range("M150").Value = 350
$ws.Range("N150").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O150").Value = "Región de O'Higgins"
$ws.Range("P150").Value = 350
$ws.Range("Q150").Value = 1
$ws.Range("R150").Value = "Hortaliza"
